$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 935 (shifts existing rows 935:967 down to 936:968,
# matching the date-ordered list of Cilantro observations for
# "Vega Central Mapocho de Santiago" with one new weekly record added).
$ws.Rows("935:935").Insert()

# Populate the newly inserted row 935 with the new observation.
$ws.Range("A935").Value = 9
$ws.Range("B935").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C935").Value = "Metropolitana"
$ws.Range("D935").Value = 45075
$ws.Range("E935").Value = 13
$ws.Range("F935").Value = 100112040
$ws.Range("G935").Value = "Cilantro"
$ws.Range("H935").Value = "Sin especificar"
$ws.Range("I935").Value = "Primera"
$ws.Range("J935").Value = 70
$ws.Range("K935").Value = 15000
$ws.Range("L935").Value = 17000
$ws.Range("M935").Value = 16000
$ws.Range("N935").Value = "$/docena de atados"
$ws.Range("O935").Value = "Región Metropolitana"
$ws.Range("P935").Value = 5333
$ws.Range("Q935").Value = 3
$ws.Range("R935").Value = "Hortaliza"

# Match the date-cell number format used by the rest of column D.
$ws.Range("D935").NumberFormat = $ws.Range("D936").NumberFormat
